# Natmi following Dr Hou advice
# Adds a new "ECs" sending/target cluster to the Ccl4-Ccr1 LR-pair sheet:
# the previously single data row (M2 -> M2) becomes four rows covering all
# combinations of ECs/M2 as sending and target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns A-D are text (cluster/symbol labels), E-T are numeric metrics
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @{
    2 = @("ECs","Ccl4","Ccr1","ECs",1,0.3333333333333333,335.6132,1006.8396,0.80464917790985,0.80464917790985,2,0.6666666666666666,114.5606336666667,343.681901,0.7368570786832789,0.736857078683279,38448.06085889773,346032.5477300796,0.592911442599554,0.5929114425995541)
    3 = @("ECs","Ccl4","Ccr1","M2",1,0.3333333333333333,335.6132,1006.8396,0.80464917790985,0.80464917790985,3,1,40.91135266666667,122.734058,0.263142921316721,0.2631429213167211,13730.38998478853,123573.5098630968,0.2117377353102959,0.211737735310296)
    4 = @("M2","Ccl4","Ccr1","ECs",3,1,81.47937800000001,244.438134,0.19535082209015,0.19535082209015,2,0.6666666666666666,114.5606336666667,343.681901,0.7368570786832789,0.736857078683279,9334.329174445862,84008.96257001275,0.1439456360837249,0.1439456360837249)
    5 = @("M2","Ccl4","Ccr1","M2",3,1,81.47937800000001,244.438134,0.19535082209015,0.19535082209015,3,1,40.91135266666667,122.734058,0.263142921316721,0.2631429213167211,3333.431568418642,30000.88411576777,0.05140518600642511,0.05140518600642511)
}

foreach ($r in 2..5) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $addr = "{0}{1}" -f $columns[$i], $r
        $ws.Range($addr).Value = $values[$i]
    }
}
